$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.7007990479469299
$ws.Range("B1").Value = 1.388329863548279
$ws.Range("C1").Value = 4.149904727935791
$ws.Range("D1").Value = 2.564904451370239
$ws.Range("E1").Value = 0.5634708404541016
